$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (Rang, Speler, Score, 180'ers, 100+ finishes, Totaal Score, Aantal Darts, 3-Darts Gemiddelde, Totaal, Winnaar)
$data = @(
    @(1, "Burger Peach", 13, 0, 1, 7345, 341, 64.62, 14, 1),
    @(2, "Yannick den Daggelder", 12, 0, 0, 6777, 370, 54.95, 12, 1),
    @(3, "Rocky Van Den Eeckhoudt", 10, 0, 1, 8003, 429, 55.97, 11, 0),
    @(4, "Niels van Dommelen", 10, 0, 0, 9695, 583, 49.89, 10, 0),
    @(5, "Lukas G", 4, 1, 0, 5294, 283, 56.12, 5, 0),
    @(6, "Nigel Riedel", 2, 0, 0, 2695, 155, 52.16, 2, 0),
    @(6, "Noah B", 2, 0, 0, 2647, 193, 41.15, 2, 0),
    @(8, "Sion Foulkes", 1, 0, 0, 3708, 247, 45.04, 1, 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]
    $ws.Cells.Item($row, 10).Value = $entry[9]
    $row++
}
